# Lab 04.docx: change the due date from "10/06." to "10/13."
#
# The target OOXML splits the single run "10/06." into three separate
# runs: "10/", "13", "." (the "13" replacing "06"). A plain Range.Text
# assignment / Find&Replace collapses the whole paragraph into one run,
# which does not match. Using tracked changes for the edit and then
# accepting only the recorded revisions (instead of Document-wide
# AcceptAllRevisions, which re-normalizes the whole document) produces
# exactly that run split while leaving the rest of the document intact.

$d = $word.ActiveDocument

$origTrackRevisions = $d.TrackRevisions

# Locate "10/06." without altering anything yet. Keep the Range in a
# variable: Find.Execute mutates the Range object it's called on to
# the hit location, but a fresh $d.Content access would just return a
# brand-new Range over the whole document again.
$searchRange = $d.Content
$found = $searchRange.Find.Execute("10/06.", $true, $false, $false, $false, `
                                    $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find '10/06.' in the document."
}

$dateStart = $searchRange.Start

# Turn on revision tracking so the inserted/deleted text is recorded as
# its own run instead of being silently merged into its neighbors.
$d.TrackRevisions = $true

# Split "10/06." right after the slash, and insert the new day "13"
# there: "10/" | "13" | "06."
$splitPos = $dateStart + 3
$insertionPoint = $d.Range($splitPos, $splitPos)
$insertionPoint.InsertBefore("13")

# Remove the old day "06", leaving: "10/" | "13" | "."
$oldDay = $d.Range($splitPos + 2, $splitPos + 4)
$oldDay.Delete()

# Restore the original track-changes setting and accept just the
# revisions we made (accepting revision-by-revision keeps the rest of
# the document's runs/rsids/lastRenderedPageBreak hints untouched,
# unlike AcceptAllRevisions which re-serializes the whole document).
$d.TrackRevisions = $origTrackRevisions
for ($i = $d.Revisions.Count; $i -ge 1; $i--) {
    $d.Revisions.Item($i).Accept()
}
